# Fix figure 2 linetype error: the "Date" column (F) was using July 1st
# of each year instead of January 1st. Recompute column F from the
# "Year" value already stored in column B for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $year = $ws.Cells.Item($r, 2).Value()
    if ($year -ne $null) {
        $jan1 = Get-Date -Year ([int]$year) -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
        $ws.Cells.Item($r, 6).Value = $jan1
    }
}
